$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column J, mirroring the style of the other header cells (A1:I1)
# Copy I1's formatting (bold font, thin border, center/top alignment) onto J1
# first, then set its text so the new header matches the look of the rest of
# row 1.
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "most_common_number_of_stories"

# Fill J2:J43 with the value 1 (most common number of stories per traveller)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 10).Value = 1
}
